$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay text so values like
# "0.05090" or "  -4.17%  " are not re-interpreted as numbers and lose
# their formatting (leading/trailing zeros, padding spaces, etc.).
$ws.Range("D2:E51").NumberFormat = "@"

# Rows 33 and 34: Filecoin and ARBITRUM swap rank order, with refreshed price/volume data
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.305"
$ws.Range("E33").Value = "  -5.00%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.105"
$ws.Range("E34").Value = "  -1.44%  "

# Remaining rows: updated Price / Volume(1h) figures
$ws.Range("D2").Value = "26.338.48"
$ws.Range("E2").Value = "  -4.17%  "
$ws.Range("D3").Value = "1.760.74"
$ws.Range("E3").Value = "  -3.50%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "304.24"
$ws.Range("E6").Value = "  -2.43%  "
$ws.Range("D7").Value = "0.4278"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("D8").Value = "0.3612"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").Value = "0.07051"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").Value = "0.8313"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").Value = "20.14"
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").Value = "1.765.47"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "5.228"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("D14").Value = "6.389"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "0.06806"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "79.17"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "0.000008634"
$ws.Range("E18").Value = "  -2.63%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "14.98"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "26.127.89"
$ws.Range("E21").Value = "  -4.22%  "
$ws.Range("D22").Value = "5.001"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("D24").Value = "1.960.69"
$ws.Range("E24").Value = "  -3.29%  "
$ws.Range("E25").Value = "  -4.15%  "
$ws.Range("D26").Value = "152.39"
$ws.Range("E26").Value = "  -1.67%  "
$ws.Range("D27").Value = "18.11"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("D28").Value = "114.83"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").Value = "5.022"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "1.684"
$ws.Range("E30").Value = "  -6.59%  "
$ws.Range("D31").Value = "0.08889"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").Value = "0.7214"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D35").Value = "2.748"
$ws.Range("E35").Value = "  -7.91%  "
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "1.069"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").Value = "0.05090"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").Value = "0.01884"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").Value = "0.4888"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").Value = "0.1599"
$ws.Range("E41").Value = "  -2.63%  "
$ws.Range("D42").Value = "6.210"
$ws.Range("E42").Value = "  -4.00%  "
$ws.Range("E43").Value = "  -10.34%  "
$ws.Range("D44").Value = "7.982"
$ws.Range("E44").Value = "  -3.65%  "
$ws.Range("D45").Value = "104.66"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "10.05"
$ws.Range("E47").Value = "  -3.53%  "
$ws.Range("D48").Value = "0.06188"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").Value = "0.4457"
$ws.Range("E49").Value = "  -4.46%  "
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").Value = "1.729"
$ws.Range("E51").Value = "  +0.55%  "
